$d = $word.ActiveDocument

# Remove the trailing "Ver no Jupiter..." and "© 2020..." paragraphs,
# together with the blank paragraph that immediately precedes them
# (the one right after the LOM3221 requisito line).

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $p.Range.Delete()
    }
    elseif ($t -like "*Contact: luizeleno@usp.br*") {
        $p.Range.Delete()
    }
}

# Now remove the now-orphaned blank paragraph that sat between the
# LOM3221 line and the (now deleted) "Ver no Jupiter" paragraph.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t.Trim() -eq "") {
        $prev = $null
        if ($i -gt 1) { $prev = $d.Paragraphs.Item($i - 1) }
        if ($prev -ne $null -and $prev.Range.Text -like "*LOM3221: Laboratório de Eletrônica (Indicação de Conjunto)*") {
            $p.Range.Delete()
            break
        }
    }
}
